# The commit swaps the two theme parts of this deck: the "Integral" theme
# (ppt/theme/theme1.xml, used by the slide master / slides) and the
# "Office Theme" (ppt/theme/theme2.xml, used only by the notes master)
# traded places - theme1.xml ends up holding the Office Theme palette and
# theme2.xml ends up holding the Integral palette.
#
# The presentation's live/addressable theme through the PowerPoint object
# model is the one tied to the slide master (theme1.xml); its 12 theme
# colors are reachable and read/write via ThemeColorScheme. Re-point them
# at the Office Theme color values so the active theme's color scheme
# becomes the Office Theme palette, matching the target theme1.xml.

$p = $ppt.ActivePresentation

# Office Theme color scheme (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
# expressed as COM RGB() values (0x00BBGGRR) computed from the hex triplets:
#   dk1=000000 lt1=FFFFFF dk2=44546A lt2=E7E6E6
#   accent1=5B9BD5 accent2=ED7D31 accent3=A5A5A5 accent4=FFC000
#   accent5=4472C4 accent6=70AD47 hlink=0563C1 folHlink=954F72
$officeThemeRgb = @(0, 16777215, 6968388, 15132391, 13998939, 3243501, 10855845, 49407, 12874308, 4697456, 12673797, 7491477)

$tcs = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Item($i).RGB = $officeThemeRgb[$i - 1]
}
